# "Add duration to multiple contacts"
#
# In the survey_3 sheet (codebook for survey_3), three new variable rows are
# inserted, one right after each of the "precautions" rows for the
# work / school / other "multiple contacts" blocks, adding a matching
# "*_duration" variable. Sheet selection/active-tab state is also updated to
# reflect where the editor was working (survey_3 becomes the active sheet).

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("survey_3")

# --- Insert the three new rows (top to bottom, so later offsets already
#     account for the rows inserted above them) -----------------------------

# After q79c / cnt_multiple_contacts_work_precautions (row 189) insert
# q79d / multiple_contacts_work_duration
$ws3.Rows.Item(190).Insert()
$ws3.Range("A190").Value = "q79d"
$ws3.Range("B190").Value = "multiple_contacts_work_duration"
$ws3.Range("A190:B190").Font.Color = 0

# After q80c / cnt_multiple_contacts_school_precautions insert
# q80d / multiple_contacts_school_duration
$ws3.Rows.Item(195).Insert()
$ws3.Range("A195").Value = "q80d"
$ws3.Range("B195").Value = "multiple_contacts_school_duration"
$ws3.Range("A195:B195").Font.Color = 0

# After q81c / cnt_multiple_contacts_other_precautions insert
# q81d / multiple_contacts_other_duration
$ws3.Rows.Item(200).Insert()
$ws3.Range("A200").Value = "q81d"
$ws3.Range("B200").Value = "multiple_contacts_other_duration"
$ws3.Range("A200:B200").Font.Color = 0

# --- Update view / selection state -----------------------------------------

# survey_4: selection moves to A172:B172
$ws4 = $wb.Worksheets.Item("survey_4")
$ws4.Range("A172:B172").Select()

# survey_5: selection moves to A159:B160 (and it stops being the active tab)
$ws5 = $wb.Worksheets.Item("survey_5")
$ws5.Range("A159:B160").Select()

# survey_3 becomes the active sheet/tab, with the new row selected
$ws3.Activate()
$ws3.Range("A200:B200").Select()
